$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 1207
$ws.Range("I115").Value = 1207
$ws.Range("K115").Value = 3621
$ws.Range("M115").Value = -2054
$ws.Range("H125").Value = 1365.375
$ws.Range("I125").Value = 314.66666
$ws.Range("J125").Value = 1995.8
$ws.Range("K125").Value = 2831.99994
$ws.Range("L125").Value = 17962.2
$ws.Range("M125").Value = -371.9999399999997
$ws.Range("N125").Value = -22882.2
$ws.Range("H136").Value = 48280
$ws.Range("J136").Value = 48280
$ws.Range("L136").Value = 48280
$ws.Range("N136").Value = -58480
$ws.Range("H137").Value = 1222952.8
$ws.Range("I137").Value = 1537473.1
$ws.Range("J137").Value = 4186.25
$ws.Range("K137").Value = 4612419.300000001
$ws.Range("L137").Value = 12558.75
$ws.Range("M137").Value = -4609869.300000001
$ws.Range("N137").Value = -17658.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2037.4546
$ws.Range("I61").Value = 2022.4
$ws.Range("J61").Value = 2050
$ws.Range("K61").Value = 2022.4
$ws.Range("L61").Value = 2050
$ws.Range("M61").Value = -1810.4
$ws.Range("N61").Value = -2474
$ws.Range("H74").Value = 2029.762
$ws.Range("I74").Value = 1157.75
$ws.Range("J74").Value = 3192.4443
$ws.Range("K74").Value = 1157.75
$ws.Range("L74").Value = 3192.4443
$ws.Range("M74").Value = -283.75
$ws.Range("N74").Value = -4940.4443
$ws.Range("H77").Value = 2029.762
$ws.Range("I77").Value = 1157.75
$ws.Range("J77").Value = 3192.4443
$ws.Range("K77").Value = 5788.75
$ws.Range("L77").Value = 15962.2215
$ws.Range("M77").Value = -1420.75
$ws.Range("N77").Value = -24698.2215
$ws.Range("H97").Value = 1807.9166
$ws.Range("I97").Value = 1019.2222
$ws.Range("J97").Value = 4174
$ws.Range("K97").Value = 1019.2222
$ws.Range("L97").Value = 4174
$ws.Range("M97").Value = -523.2222
$ws.Range("N97").Value = -5166
$ws.Range("H136").Value = 2037.4546
$ws.Range("I136").Value = 2022.4
$ws.Range("J136").Value = 2050
$ws.Range("K136").Value = 6067.200000000001
$ws.Range("L136").Value = 6150
$ws.Range("M136").Value = -3517.200000000001
$ws.Range("N136").Value = -11250

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2097
$ws.Range("I107").Value = 1829.65
$ws.Range("K107").Value = 1829.65
$ws.Range("M107").Value = 90.34999999999991
$ws.Range("H134").Value = 2295.7646
$ws.Range("I134").Value = 1467.7142
$ws.Range("J134").Value = 6160
$ws.Range("K134").Value = 4403.142599999999
$ws.Range("L134").Value = 18480
$ws.Range("M134").Value = -1868.142599999999
$ws.Range("N134").Value = -23550

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1875.75
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 2001
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 2001
$ws.Range("M16").Value = -1213
$ws.Range("N16").Value = -2575
$ws.Range("H31").Value = 5171.931
$ws.Range("I31").Value = 1084.5217
$ws.Range("K31").Value = 1084.5217
$ws.Range("M31").Value = -789.5217
$ws.Range("H34").Value = 5171.931
$ws.Range("I34").Value = 1084.5217
$ws.Range("K34").Value = 1084.5217
$ws.Range("M34").Value = -882.5217
$ws.Range("H58").Value = 2268.5522
$ws.Range("I58").Value = 1763.86
$ws.Range("J58").Value = 3752.9412
$ws.Range("K58").Value = 1763.86
$ws.Range("L58").Value = 3752.9412
$ws.Range("M58").Value = -1560.86
$ws.Range("N58").Value = -4158.9412
$ws.Range("H99").Value = 10530601
$ws.Range("I99").Value = 20002142
$ws.Range("K99").Value = 20002142
$ws.Range("M99").Value = -20000644
$ws.Range("H107").Value = 760
$ws.Range("J107").Value = 965
$ws.Range("L107").Value = 965
$ws.Range("N107").Value = -4805
$ws.Range("H113").Value = 1875.75
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2001
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2001
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6341
$ws.Range("H122").Value = 4342.4
$ws.Range("I122").Value = 1678
$ws.Range("K122").Value = 5034
$ws.Range("M122").Value = -2584
$ws.Range("H126").Value = 10530601
$ws.Range("I126").Value = 20002142
$ws.Range("K126").Value = 60006426
$ws.Range("M126").Value = -60003956
$ws.Range("H132").Value = 5484.4375
$ws.Range("I132").Value = 5159.364
$ws.Range("J132").Value = 6199.6
$ws.Range("K132").Value = 15478.092
$ws.Range("L132").Value = 18598.8
$ws.Range("M132").Value = -12948.092
$ws.Range("N132").Value = -23658.8
$ws.Range("H134").Value = 5602.3213
$ws.Range("I134").Value = 7931.533
$ws.Range("J134").Value = 2914.7693
$ws.Range("K134").Value = 23794.599
$ws.Range("L134").Value = 8744.3079
$ws.Range("M134").Value = -21259.599
$ws.Range("N134").Value = -13814.3079
$ws.Range("H136").Value = 2268.5522
$ws.Range("I136").Value = 1763.86
$ws.Range("J136").Value = 3752.9412
$ws.Range("K136").Value = 5291.58
$ws.Range("L136").Value = 11258.8236
$ws.Range("M136").Value = -2741.58
$ws.Range("N136").Value = -16358.8236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3484
$ws.Range("I126").Value = 2916.6667
$ws.Range("K126").Value = 8750.000100000001
$ws.Range("M126").Value = -6280.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6149.8335
$ws.Range("I40").Value = 4379.8
$ws.Range("K40").Value = 4379.8
$ws.Range("M40").Value = -4243.8
$ws.Range("H132").Value = 5590.8486
$ws.Range("I132").Value = 3031.5789
$ws.Range("J132").Value = 9064.143
$ws.Range("K132").Value = 9094.736699999999
$ws.Range("L132").Value = 27192.429
$ws.Range("M132").Value = -6564.736699999999
$ws.Range("N132").Value = -32252.429
$ws.Range("H136").Value = 4863.4736
$ws.Range("I136").Value = 1510.6
$ws.Range("K136").Value = 4531.799999999999
$ws.Range("M136").Value = -1981.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15158556
$ws.Range("I132").Value = 9171.154
$ws.Range("J132").Value = 37041000
$ws.Range("K132").Value = 27513.462
$ws.Range("L132").Value = 111123000
$ws.Range("M132").Value = -24983.462
$ws.Range("N132").Value = -111128060
$ws.Range("H136").Value = 6713.5586
$ws.Range("I136").Value = 6011.1304
$ws.Range("J136").Value = 8182.273
$ws.Range("K136").Value = 18033.3912
$ws.Range("L136").Value = 24546.819
$ws.Range("M136").Value = -15483.3912
$ws.Range("N136").Value = -29646.819
